# Línea 141 schedule refresh (scrape @ 20:15:00) applied via Excel COM interop.
# Three worksheets share the same shape: A=Hora_Scrap, B=Hora_Llegada, C=Linea, D=Minutos, E=Parada,
# with A2 = "Última actualización: HH:MM:SS" and A3 = "Total filas: N" as free-text headers.
$wb = $excel.ActiveWorkbook

### Sheet "LP1912" ######################################################
$ws = $wb.Worksheets.Item("LP1912")

# Header banner: new scrape timestamp + updated row count (510 -> 516 data rows).
$ws.Cells.Item(2,1).Value = "Última actualización: 20:15:00"
$ws.Cells.Item(3,1).Value = "Total filas: 516"

# A later re-scrape reordered a handful of same-arrival-time rows (tie-break now by
# Hora_Scrap ascending); only columns A/C/D move; B (Hora_Llegada) and E (Parada) are stable.
$ws.Cells.Item(117,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(118,3).Value = "16_SANTA ANA"
$ws.Cells.Item(229,1).Value = "11:48:04"
$ws.Cells.Item(229,3).Value = "10_OLMOS"
$ws.Cells.Item(229,4).Value = 92
$ws.Cells.Item(230,1).Value = "13:19:56"
$ws.Cells.Item(230,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(230,4).Value = 1
$ws.Cells.Item(245,1).Value = "13:19:56"
$ws.Cells.Item(245,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(245,4).Value = 31
$ws.Cells.Item(246,1).Value = "12:01:50"
$ws.Cells.Item(246,3).Value = "215A_EL PATO"
$ws.Cells.Item(246,4).Value = 109
$ws.Cells.Item(278,1).Value = "14:58:43"
$ws.Cells.Item(278,3).Value = "16_SANTA ANA"
$ws.Cells.Item(278,4).Value = 0
$ws.Cells.Item(279,1).Value = "13:19:56"
$ws.Cells.Item(279,3).Value = "215B_EL PATO"
$ws.Cells.Item(279,4).Value = 99
$ws.Cells.Item(310,1).Value = "14:19:48"
$ws.Cells.Item(310,3).Value = "27_EL RETIRO"
$ws.Cells.Item(310,4).Value = 97
$ws.Cells.Item(311,1).Value = "15:31:33"
$ws.Cells.Item(311,3).Value = "17_ROMERO"
$ws.Cells.Item(311,4).Value = 25
$ws.Cells.Item(375,1).Value = "16:53:01"
$ws.Cells.Item(375,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(375,4).Value = 42
$ws.Cells.Item(376,1).Value = "15:57:48"
$ws.Cells.Item(376,3).Value = "27_EL RETIRO"
$ws.Cells.Item(376,4).Value = 98
$ws.Cells.Item(407,1).Value = "18:17:05"
$ws.Cells.Item(407,3).Value = "16_SANTA ANA"
$ws.Cells.Item(407,4).Value = 1
$ws.Cells.Item(409,1).Value = "17:59:03"
$ws.Cells.Item(409,3).Value = "15_ABASTO"
$ws.Cells.Item(409,4).Value = 19
$ws.Cells.Item(440,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(441,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(455,1).Value = "18:37:25"
$ws.Cells.Item(455,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(455,4).Value = 45
$ws.Cells.Item(456,1).Value = "18:51:07"
$ws.Cells.Item(456,3).Value = "16_SANTA ANA"
$ws.Cells.Item(456,4).Value = 31
$ws.Cells.Item(457,3).Value = "14_ABASTO"
$ws.Cells.Item(470,1).Value = "18:17:05"
$ws.Cells.Item(470,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(470,4).Value = 94
$ws.Cells.Item(471,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(472,1).Value = "18:37:25"
$ws.Cells.Item(472,3).Value = "11X44_ETCHEVERRY"
$ws.Cells.Item(472,4).Value = 74
$ws.Cells.Item(479,1).Value = "19:56:21"
$ws.Cells.Item(479,3).Value = "14_ABASTO"
$ws.Cells.Item(479,4).Value = 4
$ws.Cells.Item(480,1).Value = "18:37:25"
$ws.Cells.Item(480,3).Value = "17_ROMERO"
$ws.Cells.Item(480,4).Value = 83

# Tail of the table (rows 488-521): re-sorted plus 6 brand-new rows scraped at 20:15:00,
# so every row in this block is rewritten in full.
$tailRows = @(
    @(488, "20:15:00", "20:16", "15_ABASTO", 1, "LP1912"),
    @(489, "20:15:00", "20:16", "23_HERNANDEZ", 1, "LP1912"),
    @(490, "20:15:00", "20:16", "16_SANTA ANA", 1, "LP1912"),
    @(491, "18:58:44", "20:21", "26_HERNANDEZ", 83, "LP1912"),
    @(492, "18:37:25", "20:22", "26_HERNANDEZ", 105, "LP1912"),
    @(493, "19:42:02", "20:22", "11_ETCHEVERRY", 40, "LP1912"),
    @(494, "18:37:25", "20:23", "11_ETCHEVERRY", 106, "LP1912"),
    @(495, "19:56:21", "20:23", "215A_EL PATO", 27, "LP1912"),
    @(496, "18:37:25", "20:24", "215A_EL PATO", 107, "LP1912"),
    @(497, "18:51:07", "20:25", "215A_EL PATO", 94, "LP1912"),
    @(498, "19:42:02", "20:26", "14_ABASTO", 44, "LP1912"),
    @(499, "18:51:07", "20:27", "14_ABASTO", 96, "LP1912"),
    @(500, "19:42:02", "20:31", "225_GOMEZ", 49, "LP1912"),
    @(501, "18:37:25", "20:32", "225_GOMEZ", 115, "LP1912"),
    @(502, "18:58:44", "20:35", "14_ABASTO", 97, "LP1912"),
    @(503, "19:42:02", "20:39", "11_ETCHEVERRY", 57, "LP1912"),
    @(504, "18:51:07", "20:46", "14X44_ABASTO", 115, "LP1912"),
    @(505, "18:58:44", "20:48", "14X44_ABASTO", 110, "LP1912"),
    @(506, "19:42:02", "20:52", "15_ABASTO", 70, "LP1912"),
    @(507, "19:56:21", "20:52", "23_HERNANDEZ", 56, "LP1912"),
    @(508, "19:42:02", "20:53", "23_HERNANDEZ", 71, "LP1912"),
    @(509, "18:58:44", "20:56", "10_OLMOS", 118, "LP1912"),
    @(510, "19:42:02", "20:57", "27_EL RETIRO", 75, "LP1912"),
    @(511, "19:56:21", "21:00", "215B_EL PATO", 64, "LP1912"),
    @(512, "19:42:02", "21:01", "215B_EL PATO", 79, "LP1912"),
    @(513, "19:42:02", "21:04", "84_COLONIA URQUIZA-ESC 49", 82, "LP1912"),
    @(514, "20:15:00", "21:16", "84_COLONIA URQUIZA-ESC 49", 61, "LP1912"),
    @(515, "19:42:02", "21:21", "26_HERNANDEZ", 99, "LP1912"),
    @(516, "19:42:02", "21:23", "10_OLMOS", 101, "LP1912"),
    @(517, "19:42:02", "21:38", "14_ABASTO", 116, "LP1912"),
    @(518, "19:42:02", "21:38", "17_ROMERO", 116, "LP1912"),
    @(519, "20:15:00", "21:43", "17_ROMERO", 88, "LP1912"),
    @(520, "19:56:21", "21:47", "215A_EL PATO", 111, "LP1912"),
    @(521, "20:15:00", "22:08", "17_ROMERO", 113, "LP1912")
)
foreach ($row in $tailRows) {
    $r = $row[0]
    $ws.Cells.Item($r,1).Value = $row[1]
    $ws.Cells.Item($r,2).Value = $row[2]
    $ws.Cells.Item($r,3).Value = $row[3]
    $ws.Cells.Item($r,4).Value = $row[4]
    $ws.Cells.Item($r,5).Value = $row[5]
}

### Sheet "LP1912-215" ##################################################
$ws = $wb.Worksheets.Item("LP1912-215")
# Only the "last updated" banner changes here; the 215-filtered rows themselves are untouched.
$ws.Cells.Item(2,1).Value = "Última actualización: 20:15:00"

### Sheet "6203-6173" ###################################################
$ws = $wb.Worksheets.Item("6203-6173")
# Header banner: new scrape timestamp + updated row count (67 -> 68 data rows).
$ws.Cells.Item(2,1).Value = "Última actualización: 20:15:00"
$ws.Cells.Item(3,1).Value = "Total filas: 68"

# One brand-new row appended (215A_LA PLATA arriving at 22:05).
$ws.Cells.Item(73,1).Value = "20:15:00"
$ws.Cells.Item(73,2).Value = "22:05"
$ws.Cells.Item(73,3).Value = "215A_LA PLATA"
$ws.Cells.Item(73,4).Value = 110
$ws.Cells.Item(73,5).Value = "L6173"

